# Add a "Save" column (H) to the s_vals worksheet, matching column F
# (0/1 flag) -- this is the "add save column in s_vals sheets" change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: text "Save", formatted like the other header cells
# (bold, bordered, centered) by copying the format from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for H2:H53, as specified by the diff.
$saveValues = @(
    0,0,0,0,0,0,0,0,1,1,
    0,0,0,1,1,0,0,0,0,0,
    0,0,0,0,0,1,0,0,0,0,
    1,0,0,0,0,0,0,0,0,0,
    0,0,0,1,1,0,0,0,0,0,
    0,0
)

$firstRow = 2
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 8).Value = $saveValues[$i]
}
